$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Cars_Are_Fun"
$ws.Range("B6").Value = "Cars are fun"
$ws.Range("C6").Value = "Bilar är kul"
$ws.Range("D6").Value = "Need review"

# Mirror the existing pattern in the sheet where column E has an empty,
# unstyled placeholder cell on every data row (e.g. E2:E5). Copying an
# existing empty E cell preserves that exact "no value / no style" shape
# for the new row instead of leaving the cell absent or style-stamped.
$ws.Range("E5").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$excel.CutCopyMode = 0
